# Fix a few typos across the deck (per commit message "fixed a few typos").

$p = $ppt.ActivePresentation

# --- Slide 9: Title "Using Result<T,E> with is_ok()" -> "Using Result<T, E> with is_ok()" ---
# Only the first run ("Using Result<T,E> with ") changes; later runs ("is_ok", "()") stay intact.
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(1)
$tr9 = $sh9.TextFrame.TextRange
$run9 = $tr9.Runs(1, 1)
$run9.Text = "Using Result<T, E> with "

# --- Slide 10: "If you use code..." -> "If you need to use code..." ---
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(2)
$tr10 = $sh10.TextFrame.TextRange
$para10 = $tr10.Paragraphs(1, 1)
$run10 = $para10.Runs(1, 1)
$run10.Text = "If you need to use code that doesn’t reliably avoid panics you may attempt to trap them:"

# --- Slide 13: "if let uses matching operator =" -> "“if let” uses matching operator =" ---
$s13 = $p.Slides.Item(13)
$sh13 = $s13.Shapes.Item(2)
$tr13 = $sh13.TextFrame.TextRange
$para13 = $tr13.Paragraphs(4, 1)
$run13 = $para13.Runs(1, 1)
$run13.Text = "“if let” uses matching operator ="

# --- Slide 14: "Demonstration code using match and let if" -> "...match and if let" ---
$s14 = $p.Slides.Item(14)
$sh14 = $s14.Shapes.Item(1)
$tr14 = $sh14.TextFrame.TextRange
$para14 = $tr14.Paragraphs(1, 1)
$run14 = $para14.Runs(1, 1)
$run14.Text = "Demonstration code using match and if let"

# --- Slide 15: both "Fn" runs (start of two code paragraphs) -> "fn" ---
$s15 = $p.Slides.Item(15)
$sh15 = $s15.Shapes.Item(2)
$tr15 = $sh15.TextFrame.TextRange

$para15a = $tr15.Paragraphs(1, 1)
$chars15a = $tr15.Characters($para15a.Start, 2)
$chars15a.Text = "fn"

$para15b = $tr15.Paragraphs(2, 1)
$chars15b = $tr15.Characters($para15b.Start, 2)
$chars15b.Text = "fn"

# --- Slide 22: "Fn" (with_options line, lvl 1) -> "fn" ---
$s22 = $p.Slides.Item(22)
$sh22 = $s22.Shapes.Item(2)
$tr22 = $sh22.TextFrame.TextRange
$para22 = $tr22.Paragraphs(5, 1)
$chars22 = $tr22.Characters($para22.Start, 2)
$chars22.Text = "fn"
